$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows, per repulled/mean-calculated data
$ws.Range("F3").Value = -1
$ws.Range("F5").Value = -5
$ws.Range("F11").Value = -3
$ws.Range("F21").Value = 4
